# InvestmentCalc.xlsx update
# - Extend the projection from 10 years (cols B:L) to 15 years (cols B:Q)
# - Re-parameterise Depreciation / Incoming / Outgoing payments
# - Move the one-off Residual / restricted-Equity recovery from the old
#   last year (L) to the new last year (Q)
# - Recompute Yearly Net / Present Value / Accumulated Present Value / NPV

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helpers: column letters for year columns 0..15 (B..Q) and a small
# style-clone helper built on Copy + PasteSpecial(xlPasteFormats) so new
# cells pick up the exact same cell style as their row neighbours.
# ---------------------------------------------------------------------
$yearCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
$xlPasteFormats = -4122

function Clone-Style($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# Row 1 - Year header: add years 11..15 in M1:Q1 (same style as L1)
# ---------------------------------------------------------------------
Clone-Style "L1" "M1:Q1"
for ($i = 11; $i -le 15; $i++) {
    $col = $yearCols[$i]
    $ws.Range($col + "1").Value = $i
}

# ---------------------------------------------------------------------
# Row 3 - Depreciation: 30000 -> 20000 for every payment year (C:Q)
# ---------------------------------------------------------------------
Clone-Style "L3" "M3:Q3"
for ($i = 1; $i -le 15; $i++) {
    $col = $yearCols[$i]
    $ws.Range($col + "3").Value = 20000
}

# ---------------------------------------------------------------------
# Row 4 - Incoming Payments: 1400000 -> 1050000 for every payment year (C:Q)
# ---------------------------------------------------------------------
Clone-Style "L4" "M4:Q4"
for ($i = 1; $i -le 15; $i++) {
    $col = $yearCols[$i]
    $ws.Range($col + "4").Value = 1050000
}

# ---------------------------------------------------------------------
# Row 5 - Outgoing Payments: year0 -350000 -> -70000, years1..15 -700000 -> -350000
# ---------------------------------------------------------------------
Clone-Style "L5" "M5:Q5"
$ws.Range("B5").Value = -70000
for ($i = 1; $i -le 15; $i++) {
    $col = $yearCols[$i]
    $ws.Range($col + "5").Value = -350000
}

# ---------------------------------------------------------------------
# Row 6 - Residual: moves from L6 (140000) to Q6 (70000, half of before)
# Row 7 - restricted Equity recovery: moves from L7 (200000) to Q7 (200000)
# ---------------------------------------------------------------------
$ws.Range("L6").Clear()
Clone-Style "L3" "Q6"
$ws.Range("Q6").Value = 70000

$ws.Range("L7").Clear()
Clone-Style "L3" "Q7"
$ws.Range("Q7").Value = 200000

# ---------------------------------------------------------------------
# Row 8 - Yearly Net = Initial Investment + Depreciation + Incoming +
#         Outgoing + Residual + restricted Equity (per column)
# ---------------------------------------------------------------------
$initialInvestment = -1000000
$restrictedEquity = -200000

for ($i = 0; $i -le 15; $i++) {
    $col = $yearCols[$i]

    $dep = 0
    if ($i -ge 1) { $dep = $ws.Range($col + "3").Value2 }

    $incoming = 0
    if ($i -ge 1) { $incoming = $ws.Range($col + "4").Value2 }

    $outgoing = $ws.Range($col + "5").Value2

    $residual = 0
    if ($i -eq 15) { $residual = $ws.Range($col + "6").Value2 }

    $restricted = 0
    if ($i -eq 0) { $restricted = $restrictedEquity }
    if ($i -eq 15) { $restricted = $ws.Range($col + "7").Value2 }

    $initial = 0
    if ($i -eq 0) { $initial = $initialInvestment }

    $yearlyNet = $initial + $dep + $incoming + $outgoing + $residual + $restricted

    if ($i -ge 11) {
        # new columns M..Q need the style cloned first (same style as L8)
        Clone-Style "L8" ($col + "8")
    }
    $ws.Range($col + "8").Value = $yearlyNet
}

# ---------------------------------------------------------------------
# Row 9 - Present Value = Yearly Net / (1 + afterTaxDiscountRate) ^ year
# Row 10 - Accumulated Present Value = running sum of Present Value,
#          coloured green (style of L9, fillId 3) when >= 0 and red
#          (style of B9, fillId 2) when negative
# ---------------------------------------------------------------------
$discountRate = $ws.Range("B14").Value2

$accumulated = 0
for ($i = 0; $i -le 15; $i++) {
    $col = $yearCols[$i]
    $yearlyNet = $ws.Range($col + "8").Value2
    $pv = $yearlyNet / [Math]::Pow((1 + $discountRate), $i)

    if ($i -ge 11) {
        Clone-Style "L9" ($col + "9")
    }
    $ws.Range($col + "9").Value = $pv

    $accumulated = $accumulated + $pv

    if ($i -ge 11) {
        Clone-Style "L10" ($col + "10")
    }
    if ($i -gt 0) {
        if ($accumulated -lt 0) {
            Clone-Style "B10" ($col + "10")
        } else {
            Clone-Style "L9" ($col + "10")
        }
    }
    $ws.Range($col + "10").Value = $accumulated
}

# ---------------------------------------------------------------------
# Row 11 - Net Present Value = last Accumulated Present Value (Q10)
# ---------------------------------------------------------------------
$ws.Range("B11").Value = $ws.Range("Q10").Value2
